# Requests Log Form - F4.xlsx : "sops Update 4"
#
# Re-labels the form from the Change-Request (CR) series to the Software
# Development (SD) series: renames the visible worksheet tab, updates the
# revision/date stamp baked into the footer, and repoints the saved
# workbook window / last-known folder bookkeeping to the new SOP folder.

$wb = $excel.ActiveWorkbook

# --- Rename the visible worksheet: F-SW-CR-02 -> F-SW-SD-04 -----------------
$ws = $wb.Sheets("F-SW-CR-02")
$ws.Name = "F-SW-SD-04"

# --- Footer: "Rev : 0 (0/0/2025)" -> "Rev:0(01/10/2025)" --------------------
# The odd footer is built from 3 sections (&L left / &C center / &R right);
# only the right-hand "Rev" section changes.
$ws.PageSetup.RightFooter = "&14Rev:0(01/10/2025)"

# --- Make sure the renamed sheet stays the active/selected tab --------------
$ws.Select()

# --- Best-effort window/view bookkeeping ------------------------------------
# Mirrors the saved view moving from Page Break Preview (scrolled to A3) to
# Page Layout view (scrolled to A10), and the workbook window being
# maximized/repositioned on the author's second monitor. Wrapped defensively
# so the core rename/footer edits above always land even if a given property
# isn't available in this host.
try {
    $excel.ActiveWindow.View = 3            # xlPageLayoutView (was xlPageBreakPreview)
    $excel.ActiveWindow.ScrollRow = 10      # top-left row -> A10 (was A3)
    $excel.ActiveWindow.ScrollColumn = 1
    $excel.ActiveWindow.WindowState = -4143 # xlNormal
    $excel.ActiveWindow.Left = 20370
    $excel.ActiveWindow.Top = -2595
    $excel.ActiveWindow.Width = 29040
    $excel.ActiveWindow.Height = 15840
} catch {
    Write-Output "window/view properties not available: $_"
}

Write-Output "Renamed sheet to '$($ws.Name)'; footer now '$($ws.PageSetup.LeftFooter)$($ws.PageSetup.CenterFooter)$($ws.PageSetup.RightFooter)'"
